# Generate Report for Handback
# Updates the "395edba5-4b17-4623-8a01-c09bbe204788" row (row 7) on the
# zh-cn and de-de sheets with the newly generated handback info: the
# handback file hyperlink, the handback target file name, the handback
# datetime, and an error detail noting the handback is stale. Also widens
# the "Error Detail" column (P) so the new message is readable.

$wb = $excel.ActiveWorkbook

$current_sha = "a5e65809b299abee40d866c355cc0cb3e6650ae4"
$latest_sha  = "76b5ab6c87f302b58eea6d0dc3bed49a544b03b0"
$md_name     = "395edba5-4b17-4623-8a01-c09bbe204788.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$current_sha/e2e/$md_name, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$latest_sha/e2e/$md_name."
$currentUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$current_sha/e2e/$md_name"

# width 40 in the saved XML == ColumnWidth 40 - 5/6 (Excel's character-width
# padding) when expressed through the COM ColumnWidth property.
$colWidthForStored40 = 40 - (5 / 6)

$locales = @(
    @{ Sheet = "zh-cn"; Xlf = "395edba5-4b17-4623-8a01-c09bbe204788.2de445ee30f54352b90657c9d22147147894cea9.zh-cn.xlf"; Datetime = "2016-08-22 10:23:32" },
    @{ Sheet = "de-de"; Xlf = "395edba5-4b17-4623-8a01-c09bbe204788.2de445ee30f54352b90657c9d22147147894cea9.de-de.xlf"; Datetime = "2016-08-22 10:23:39" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # I7: Latest Handback File -> new hyperlink to the handback .md file
    $ws.Hyperlinks.Add($ws.Range("I7"), $currentUrl, "", "", $md_name)

    # J7: Latest Handback File -> handback xlf target file name
    $ws.Range("J7").Value = $locale.Xlf

    # K7: Latest Handback DateTime -> timestamp the handback was generated
    $ws.Range("K7").Value = $locale.Datetime

    # P7: Error Detail -> stale-handback warning
    $ws.Range("P7").Value = $errorDetail

    # Widen the Error Detail column so the message is readable
    $ws.Columns.Item(16).ColumnWidth = $colWidthForStored40
}
